$p = $ppt.ActivePresentation
$nm = $p.NotesMaster
$d = $nm.Design
Write-Output "Design name=[$($d.Name)]"
Write-Output ($d | Get-Member | Out-String)
